$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "Through 2022-06-25" to "Through 2022-06-26"
$ws.Name = "Through 2022-06-26"

# Update the column header (shared string) for the "2022 (through 06-25)" column
$ws.Range("I1").Value = "2022 (through 06-26)"

# Update June 2022 value (row 7) from 118 to 125
$ws.Range("I7").Value = 125

# Update Total 2022 value (row 14) from 781 to 788
$ws.Range("I14").Value = 788
